$d = $word.ActiveDocument

$replacements = @(
    @("741×8=5928", "605×4=2420"),
    @("932×3=2796", "573×9=5157"),
    @("499×3=1497", "241×9=2169"),
    @("226×6=1356", "729×5=3645"),
    @("576×7=4032", "953×7=6671"),
    @("495×6=2970", "364×7=2548"),
    @("769×7=5383", "217×8=1736"),
    @("571×8=4568", "578×2=1156"),
    @("197×3=591",  "190×6=1140"),
    @("395×2=790",  "952×9=8568"),
    @("407×5=2035", "460×3=1380"),
    @("233×6=1398", "786×7=5502"),
    @("992×2=1984", "163×8=1304"),
    @("309×4=1236", "556×4=2224"),
    @("336×6=2016", "214×4=856"),
    @("962×8=7696", "353×6=2118"),
    @("963×2=1926", "769×6=4614"),
    @("189×5=945",  "365×8=2920"),
    @("191×3=573",  "309×6=1854"),
    @("640×4=2560", "660×5=3300"),
    @("109×6=654",  "890×6=5340"),
    @("571×3=1713", "554×4=2216"),
    @("920×8=7360", "379×7=2653"),
    @("784×6=4704", "498×2=996"),
    @("190×9=1710", "112×3=336")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
